$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("B2").Value = 0.1843137254901961
$ws.Range("C2").Value = 0.5686274509803921
$ws.Range("P2").Value = 0.1490196078431373
$ws.Range("S2").Value = 0.09803921568627451
$ws.Range("B3").Value = 0.01342281879194631
$ws.Range("C3").Value = 0.03355704697986577
$ws.Range("J3").Value = 0.02684563758389262
$ws.Range("P3").Value = 0.7046979865771812
$ws.Range("S3").Value = 0.2214765100671141
$ws.Range("J4").Value = 0.02127659574468085
$ws.Range("P4").Value = 0.6382978723404256
$ws.Range("S4").Value = 0.3404255319148936
$ws.Range("B6").Value = 0.05612244897959184
$ws.Range("D6").Value = 0.03061224489795918
$ws.Range("F6").Value = 0.08163265306122448
$ws.Range("J6").Value = 0.1938775510204082
$ws.Range("O6").Value = 0.00510204081632653
$ws.Range("Q6").Value = 0.163265306122449
$ws.Range("R6").Value = 0.0663265306122449
$ws.Range("S6").Value = 0.4030612244897959
$ws.Range("B7").Value = 0.09644670050761421
$ws.Range("D7").Value = 0.02538071065989848
$ws.Range("E7").Value = 0.005076142131979695
$ws.Range("F7").Value = 0.05583756345177665
$ws.Range("J7").Value = 0.1573604060913706
$ws.Range("O7").Value = 0.01015228426395939
$ws.Range("Q7").Value = 0.233502538071066
$ws.Range("R7").Value = 0.05076142131979695
$ws.Range("S7").Value = 0.3654822335025381
$ws.Range("B8").Value = 0.105708245243129
$ws.Range("D8").Value = 0.01691331923890063
$ws.Range("F8").Value = 0.06765327695560254
$ws.Range("J8").Value = 0.1014799154334038
$ws.Range("O8").Value = 0.01902748414376321
$ws.Range("Q8").Value = 0.1881606765327695
$ws.Range("R8").Value = 0.06553911205073996
$ws.Range("S8").Value = 0.4355179704016913
$ws.Range("B9").Value = 0.1095890410958904
$ws.Range("D9").Value = 0.02054794520547945
$ws.Range("F9").Value = 0.06164383561643835
$ws.Range("J9").Value = 0.0958904109589041
$ws.Range("O9").Value = 0.03424657534246575
$ws.Range("Q9").Value = 0.2465753424657534
$ws.Range("R9").Value = 0.06164383561643835
$ws.Range("S9").Value = 0.3698630136986301
$ws.Range("B10").Value = 0.1011787819253438
$ws.Range("D10").Value = 0.02455795677799607
$ws.Range("F10").Value = 0.0756385068762279
$ws.Range("J10").Value = 0.1070726915520629
$ws.Range("O10").Value = 0.01571709233791748
$ws.Range("Q10").Value = 0.206286836935167
$ws.Range("R10").Value = 0.07760314341846758
$ws.Range("S10").Value = 0.3919449901768173
$ws.Range("G11").Value = 0.1140939597315436
$ws.Range("J11").Value = 0.1107382550335571
$ws.Range("K11").Value = 0.1778523489932886
$ws.Range("L11").Value = 0.587248322147651
$ws.Range("S11").Value = 0.01006711409395973
$ws.Range("G12").Value = 0.7158469945355191
$ws.Range("J12").Value = 0.1912568306010929
$ws.Range("K12").Value = 0.00546448087431694
$ws.Range("L12").Value = 0.03278688524590164
$ws.Range("S12").Value = 0.0546448087431694
$ws.Range("G13").Value = 0.7115384615384616
$ws.Range("J13").Value = 0.2692307692307692
$ws.Range("S13").Value = 0.01923076923076923
$ws.Range("F15").Value = 0.02777777777777778
$ws.Range("H15").Value = 0.1777777777777778
$ws.Range("I15").Value = 0.06111111111111111
$ws.Range("J15").Value = 0.2777777777777778
$ws.Range("K15").Value = 0.1222222222222222
$ws.Range("M15").Value = 0.01666666666666667
$ws.Range("O15").Value = 0.03333333333333333
$ws.Range("S15").Value = 0.2833333333333333
$ws.Range("F16").Value = 0.02976190476190476
$ws.Range("H16").Value = 0.1845238095238095
$ws.Range("I16").Value = 0.06547619047619048
$ws.Range("J16").Value = 0.3392857142857143
$ws.Range("K16").Value = 0.1666666666666667
$ws.Range("M16").Value = 0.01785714285714286
$ws.Range("N16").Value = 0.005952380952380952
$ws.Range("O16").Value = 0.07142857142857142
$ws.Range("S16").Value = 0.119047619047619
$ws.Range("F17").Value = 0.009779951100244499
$ws.Range("H17").Value = 0.2200488997555012
$ws.Range("I17").Value = 0.06845965770171149
$ws.Range("J17").Value = 0.4009779951100245
$ws.Range("K17").Value = 0.1075794621026895
$ws.Range("M17").Value = 0.03422982885085574
$ws.Range("N17").Value = 0.004889975550122249
$ws.Range("O17").Value = 0.05134474327628362
$ws.Range("S17").Value = 0.1026894865525672
$ws.Range("F18").Value = 0.007142857142857143
$ws.Range("H18").Value = 0.1642857142857143
$ws.Range("I18").Value = 0.1
$ws.Range("J18").Value = 0.3571428571428572
$ws.Range("K18").Value = 0.1357142857142857
$ws.Range("M18").Value = 0.04285714285714286
$ws.Range("O18").Value = 0.07857142857142857
$ws.Range("S18").Value = 0.1142857142857143
$ws.Range("F19").Value = 0.01323918799646955
$ws.Range("H19").Value = 0.2665489849955869
$ws.Range("I19").Value = 0.07149161518093557
$ws.Range("J19").Value = 0.3459841129744042
$ws.Range("K19").Value = 0.1094439541041483
$ws.Range("M19").Value = 0.02383053839364519
$ws.Range("O19").Value = 0.06266548984995587
$ws.Range("S19").Value = 0.1067961165048544
